$wb = $excel.ActiveWorkbook
$xlPasteValues = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues

# ---------------------------------------------------------------------------
# 1) "CHI TIẾT DOANH THU" (sheet #1): cyclically rotate columns A/B/C so that
#      new A = old C (Ngày thực hiện)
#      new B = old A (Tiền tố)
#      new C = old B (Mã dịch vụ)
#    Header + the 9 data rows (rows 1-10). We use Copy/PasteSpecial(values)
#    so that text that looks like a date ("07-01-2024") is moved verbatim
#    instead of being re-interpreted as a real date value.
# ---------------------------------------------------------------------------
$wsRevenue = $wb.Worksheets.Item(1)
$scratchRevenue = $wsRevenue.Cells.Item(1, 26)   # Z1 - unused scratch cell

for ($r = 1; $r -le 10; $r++) {
    $wsRevenue.Cells.Item($r, 3).Copy()
    $scratchRevenue.PasteSpecial($xlPasteValues)

    $wsRevenue.Cells.Item($r, 2).Copy()
    $wsRevenue.Cells.Item($r, 3).PasteSpecial($xlPasteValues)

    $wsRevenue.Cells.Item($r, 1).Copy()
    $wsRevenue.Cells.Item($r, 2).PasteSpecial($xlPasteValues)

    $scratchRevenue.Copy()
    $wsRevenue.Cells.Item($r, 1).PasteSpecial($xlPasteValues)
}
$scratchRevenue.Clear()
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Insert a brand-new sheet "CHI TIẾT CHI TIÊU" right after
#    "CHI TIẾT DOANH THU" and fill it with the detailed expense rows.
# ---------------------------------------------------------------------------
$wsExpenseDetail = $wb.Worksheets.Add($null, $wsRevenue)
$wsExpenseDetail.Name = "CHI TIẾT CHI TIÊU"

$wsExpenseDetail.Cells.Item(1, 1).Value = "Tiền tố"
$wsExpenseDetail.Cells.Item(1, 2).Value = "Mã chi tiêu"
$wsExpenseDetail.Cells.Item(1, 3).Value = "Ngày chi"
$wsExpenseDetail.Cells.Item(1, 4).Value = "Cơ sở"
$wsExpenseDetail.Cells.Item(1, 5).Value = "Phân loại"
$wsExpenseDetail.Cells.Item(1, 6).Value = "Lượng chi"

$scratchExpense = $wsExpenseDetail.Cells.Item(1, 26)   # Z1 - unused scratch cell

function Set-DateText($sheet, $row, $col, $text, $scratch) {
    # Write as a formula that evaluates to the literal text so it is stored
    # as a string (not auto-converted to a real date), then copy just the
    # resulting value across - this avoids creating any new cell style.
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $sheet.Cells.Item($row, $col).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
}

$expenseRows = @(
    @{ A = "CT"; B = 573; C = "07-01-2024"; D = "CẦN THƠ"; E = "Ứng Lương";          F = 200000 },
    @{ A = "CT"; B = 580; C = "07-02-2024"; D = "CẦN THƠ"; E = $null;                F = 5000000 },
    @{ A = "CT"; B = 581; C = "07-03-2024"; D = "CẦN THƠ"; E = $null;                F = $null },
    @{ A = "CT"; B = 582; C = "07-03-2024"; D = "CẦN THƠ"; E = "Tiền Thuế";          F = $null },
    @{ A = "CT"; B = 590; C = "07-03-2024"; D = "CẦN THƠ"; E = "Ứng Lương";          F = 679000 },
    @{ A = "CT"; B = 591; C = "07-03-2024"; D = "CẦN THƠ"; E = "Chi Phí Vận Hành";   F = 405000 },
    @{ A = "CT"; B = 592; C = "07-03-2024"; D = "CẦN THƠ"; E = "Chi Phí Vận Hành";   F = 1500000 },
    @{ A = "CT"; B = 599; C = "07-06-2024"; D = "CẦN THƠ"; E = "Ứng Lương";          F = 500000 }
)

$row = 2
foreach ($item in $expenseRows) {
    $wsExpenseDetail.Cells.Item($row, 1).Value = $item.A
    $wsExpenseDetail.Cells.Item($row, 2).Value = $item.B
    Set-DateText $wsExpenseDetail $row 3 $item.C $scratchExpense
    $wsExpenseDetail.Cells.Item($row, 4).Value = $item.D
    if ($item.E -ne $null) {
        $wsExpenseDetail.Cells.Item($row, 5).Value = $item.E
    }
    if ($item.F -ne $null) {
        $wsExpenseDetail.Cells.Item($row, 6).Value = $item.F
    }
    $row++
}

$scratchExpense.Clear()
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Fix the typo / rename the remaining sheets now that the new sheet has
#    shifted their positions:
#      "DAONH SỐ CÁ NHÂN" (now #3) -> "DOANH SỐ CÁ NHÂN"
#      "CHI TIÊU"          (now #4) -> "CHI TIÊU TỔNG HỢP"
#      "LŨY KẾ NGÀY"       (now #5) stays the same
# ---------------------------------------------------------------------------
$wsPersonalSales = $wb.Worksheets.Item(3)
$wsPersonalSales.Name = "DOANH SỐ CÁ NHÂN"

$wsExpenseSummary = $wb.Worksheets.Item(4)
$wsExpenseSummary.Name = "CHI TIÊU TỔNG HỢP"

# ---------------------------------------------------------------------------
# 4) Restore the original active sheet selection.
# ---------------------------------------------------------------------------
$wsFirst = $wb.Worksheets.Item(1)
$wsFirst.Activate()
